$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92
$ws.Cells.Item($row, 1).Value = 46041
$ws.Cells.Item($row, 2).Value = 214
$ws.Cells.Item($row, 3).Value = 220
$ws.Cells.Item($row, 4).Value = 207

$ws.Range("A$row").NumberFormat = $ws.Range("A91").NumberFormat
